$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for "Ngày sản xuất" / "Ngày hết hạn" (import date / expiry date)
$ws.Range("I1").Value = "Ngày sản xuất `n(dd/MM/yyyy)"
$ws.Range("J1").Value = "Ngày hết hạn`n(dd/MM/yyyy)"
